$d = $word.ActiveDocument

# 1) Append two empty paragraphs and a third paragraph containing
#    "Hello world" after the final picture, re-using the same
#    paragraph formatting (ListBullet style, no numbering, hanging
#    indent) already used at the end of the document. Do this before
#    touching proofing flags below so the new paragraph marks don't
#    inherit a stray <w:noProof/>.
$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()

$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$lastPara.Range.Text = "Hello world"

# 2) Mark every inline picture's run as "not proofed" (adds <w:noProof/>
#    to the run's rPr - either creating a new rPr or augmenting an
#    existing one), matching the 12 inline drawings in the document.
$shapeCount = $d.InlineShapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}

# 3) Merge the four separate runs that together spell out
#    git --global user.email " " into a single run.
$d.Content.Find.Execute(
    "git --global user.email " + [char]8220 + " " + [char]8220,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "git --global user.email " + [char]8220 + " " + [char]8220, 2
)
